$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# 1) Remove the whole "Gaegu" font-sample paragraph.
$d.Paragraphs.Item(5).Range.Delete()

# 2) Remove the whole "The Girl Next Door" font-sample paragraph.
#    (After step 1 it has shifted from index 12 down to index 11.)
$d.Paragraphs.Item(11).Range.Delete()

# 3) The stale <w:lastRenderedPageBreak/> cache marker needs to move off the
#    "Indie Flower" run and onto the "Nanum Pen" run (now that two
#    paragraphs earlier in the story were deleted, pagination shifted by
#    one paragraph). Rewrite both paragraphs' XML to relocate the marker.
$indieXml = "<w:p $wNs>" + `
    "<w:pPr><w:rPr>" + `
        "<w:rFonts w:ascii='Indie Flower' w:hAnsi='Indie Flower'/>" + `
        "<w:sz w:val='48'/><w:szCs w:val='48'/>" + `
    "</w:rPr></w:pPr>" + `
    "<w:r><w:rPr>" + `
        "<w:rFonts w:ascii='Indie Flower' w:hAnsi='Indie Flower'/>" + `
        "<w:sz w:val='48'/><w:szCs w:val='48'/>" + `
    "</w:rPr>" + `
    "<w:t>This is a paragraph of text so I can see how well the font looks with it.</w:t>" + `
    "</w:r></w:p>"
$d.Paragraphs.Item(7).Range.InsertXML($indieXml)

$nanumXml = "<w:p $wNs>" + `
    "<w:pPr><w:rPr>" + `
        "<w:rFonts w:ascii='Nanum Pen' w:eastAsia='Nanum Pen' w:hAnsi='Nanum Pen'/>" + `
        "<w:sz w:val='48'/><w:szCs w:val='48'/>" + `
    "</w:rPr></w:pPr>" + `
    "<w:r><w:rPr>" + `
        "<w:rFonts w:ascii='Nanum Pen' w:eastAsia='Nanum Pen' w:hAnsi='Nanum Pen'/>" + `
        "<w:sz w:val='48'/><w:szCs w:val='48'/>" + `
    "</w:rPr>" + `
    "<w:lastRenderedPageBreak/>" + `
    "<w:t>This is a paragraph of text so I can see how well the font looks with it.</w:t>" + `
    "</w:r></w:p>"
$d.Paragraphs.Item(8).Range.InsertXML($nanumXml)

# 4) Move the "_GoBack" bookmark off the trailing empty paragraph and onto
#    the start of the "Gamja Flower" paragraph.
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

$gamja = $d.Paragraphs.Item(5)
$gamjaStart = $d.Range($gamja.Range.Start, $gamja.Range.Start)
$d.Bookmarks.Add("_GoBack", $gamjaStart)
